# Fruta / hortaliza, semanal
# Update the Fecha (D), Volumen (M), Precio mínimo (N), Precio máximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns for rows 2-10
# so each row picks up the values that previously belonged to another row
# (a weekly re-shuffle of the daily records).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, for columns D, M, N, O, P, S
$data = @{
    2  = @{ D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    3  = @{ D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
    4  = @{ D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 }
    5  = @{ D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    6  = @{ D = 44517; M = 400; N = 5500;  O = 6000;  P = 5750;  S = 2875 }
    7  = @{ D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 }
    8  = @{ D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 }
    9  = @{ D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 }
    10 = @{ D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
